$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Sheet1" -> "room_id_translations"
$ws.Name = "room_id_translations"

# --- Add the new "bbc_c" block (rows 41-56) ---
# Column A is populated first for ALL new rows (bbc_c_1 .. bbc_c_16) so the
# shared-string table gets these 16 entries appended before the "BBC C0#"
# labels, matching the authoring order of the source edit.
for ($i = 1; $i -le 16; $i++) {
    $ws.Cells.Item(40 + $i, 1).Value = "bbc_c_$i"
}

# Column B only has display labels for rows 42-49 (bbc_c_2 .. bbc_c_9),
# i.e. "BBC C02".."BBC C09". Rows 41 and 50-56 have no label in column B.
for ($i = 2; $i -le 9; $i++) {
    $row = 40 + $i
    $ws.Cells.Item($row, 2).Value = "BBC C0$i"
}

# --- Update the view state to match the author's saved selection ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 85
$ws.Range("F46").Select()
